$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.784.13'
$ws.Cells.Item(2, 5).Value = '  +1.94%  '

$ws.Cells.Item(3, 4).Value = '2.114.64'
$ws.Cells.Item(3, 5).Value = '  +6.63%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  +0.08%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '333.14'

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.000'
$ws.Cells.Item(6, 5).Value = '  +0.14%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5316'
$ws.Cells.Item(7, 5).Value = '  +3.95%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4377'
$ws.Cells.Item(8, 5).Value = '  +6.57%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.08999'
$ws.Cells.Item(9, 5).Value = '  +6.64%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '45.98'
$ws.Cells.Item(10, 5).Value = '  +8.02%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.178'
$ws.Cells.Item(11, 5).Value = '  +3.91%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '25.11'
$ws.Cells.Item(12, 5).Value = '  +3.85%  '

$ws.Cells.Item(13, 4).Value = '2.107.86'
$ws.Cells.Item(13, 5).Value = '  +7.18%  '

$ws.Cells.Item(14, 5).Value = '  +4.35%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.823'
$ws.Cells.Item(15, 5).Value = '  +5.78%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '97.34'
$ws.Cells.Item(16, 5).Value = '  +3.83%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '1.001'
$ws.Cells.Item(17, 5).Value = '  +0.46%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001126'
$ws.Cells.Item(18, 5).Value = '  +1.92%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06664'
$ws.Cells.Item(19, 5).Value = '  +1.85%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '19.12'
$ws.Cells.Item(20, 5).Value = '  +1.78%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.000'
$ws.Cells.Item(21, 5).Value = '  +0.18%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.348'
$ws.Cells.Item(22, 5).Value = '  +4.45%  '

$ws.Cells.Item(23, 4).Value = '30.843.64'
$ws.Cells.Item(23, 5).Value = '  +1.93%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '12.37'
$ws.Cells.Item(24, 5).Value = '  +7.70%  '

$ws.Cells.Item(25, 4).Value = '2.356.23'
$ws.Cells.Item(25, 5).Value = '  +7.26%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.266'
$ws.Cells.Item(26, 5).Value = '  +2.76%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '22.77'
$ws.Cells.Item(27, 5).Value = '  +1.00%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.582'
$ws.Cells.Item(28, 5).Value = '  +8.69%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '162.84'
$ws.Cells.Item(29, 5).Value = '  -0.04%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '133.18'
$ws.Cells.Item(30, 5).Value = '  +1.87%  '

$ws.Cells.Item(31, 5).Value = '  +3.08%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.1079'
$ws.Cells.Item(32, 5).Value = '  +2.25%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.217'
$ws.Cells.Item(33, 5).Value = '  +3.15%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.021'
$ws.Cells.Item(34, 5).Value = '  +5.24%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.549'
$ws.Cells.Item(35, 5).Value = '  +17.88%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.02615'
$ws.Cells.Item(36, 5).Value = '  +5.69%  '

$ws.Cells.Item(37, 2).Value = 'Aptos'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '12.87'
$ws.Cells.Item(37, 5).Value = '  +9.05%  '

$ws.Cells.Item(38, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '5.527'
$ws.Cells.Item(38, 5).Value = '  +2.84%  '

$ws.Cells.Item(39, 2).Value = 'FraxShare'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '9.540'
$ws.Cells.Item(39, 5).Value = '  +7.24%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.06732'
$ws.Cells.Item(40, 5).Value = '  +3.63%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.2284'
$ws.Cells.Item(41, 5).Value = '  +5.07%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.6862'
$ws.Cells.Item(42, 5).Value = '  +4.19%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.251'
$ws.Cells.Item(43, 5).Value = '  +2.44%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.6437'
$ws.Cells.Item(44, 5).Value = '  +5.25%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.9998'
$ws.Cells.Item(45, 5).Value = '  +0.29%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '14.10'
$ws.Cells.Item(46, 5).Value = '  +4.12%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.231'
$ws.Cells.Item(47, 5).Value = '  +2.25%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '3.666'
$ws.Cells.Item(48, 5).Value = '  +0.74%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.275'
$ws.Cells.Item(49, 5).Value = '  +4.43%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '82.59'
$ws.Cells.Item(50, 5).Value = '  +3.80%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '120.42'
$ws.Cells.Item(51, 5).Value = '  -2.75%  '

